$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Generated by" label: drop the unused Admin persona, switch
#    the report to the Client persona.
$ws.Range("C3").Value = "John Snow (Client)"

# 2) The report used to show two side-by-side sections (Client info in
#    B:F, Admin info in G:I). The Admin section is no longer needed, so
#    remove columns G:I entirely (shifts nothing in since they were the
#    last columns) - this also drops the now-unused "Admin"/"Name"/
#    "Pronouns"/"Works At" header strings and collapses the header merge.
$ws.Range("G:I").Delete()

# 3) The data table had duplicated/padding rows (16-19) that are no longer
#    needed now that the sheet only has one section - remove them.
$ws.Range("16:19").Delete()

# 4) Restore the selection/active-cell bookmark saved in the sheet view.
$ws.Range("I9,G9:F9").Select()
